$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest scrape.
# Force Text format so numeric-looking strings (e.g. "27.526.77", "0.5359")
# are preserved exactly as text instead of being parsed as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.526.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.724.10"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.41"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5359"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2662"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.67"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07722"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.609"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.728.29"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.961.37"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5840"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8303"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.90"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.540.43"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +5.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.09"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +14.61%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.720"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.11%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.088"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.32"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.733"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +15.40%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.404"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.51%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.57%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.302"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.547"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.450"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.661"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.840"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9595"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.425"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5947"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01649"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.924"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8560"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.055.55"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.54%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.66"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.868.41"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.27%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "58.98"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.53%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.196"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.80%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4436"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.29%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05251"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.64%  "
